$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 634.4
$ws.Range("I2").Value = 594
$ws.Range("K2").Value = 594
$ws.Range("M2").Value = -481
$ws.Range("H18").Value = 6718.9
$ws.Range("I18").Value = 7699
$ws.Range("J18").Value = 2798.5
$ws.Range("K18").Value = 7699
$ws.Range("L18").Value = 2798.5
$ws.Range("M18").Value = -7415
$ws.Range("N18").Value = -3366.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6584.8667
$ws.Range("J45").Value = 6985.5713
$ws.Range("L45").Value = 6985.5713
$ws.Range("N45").Value = -7739.5713
$ws.Range("H110").Value = 8142.9443
$ws.Range("I110").Value = 12207.4
$ws.Range("J110").Value = 3062.375
$ws.Range("K110").Value = 12207.4
$ws.Range("L110").Value = 3062.375
$ws.Range("M110").Value = -10162.4
$ws.Range("N110").Value = -7152.375
$ws.Range("H135").Value = 75299.75
$ws.Range("J135").Value = 75299.75
$ws.Range("L135").Value = 75299.75
$ws.Range("N135").Value = -85439.75

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 41999
$ws.Range("I75").Value = 41999
$ws.Range("K75").Value = 41999
$ws.Range("M75").Value = -41063
$ws.Range("H78").Value = 41999
$ws.Range("I78").Value = 41999
$ws.Range("K78").Value = 125997
$ws.Range("M78").Value = -121317
$ws.Range("H80").Value = 1267.625
$ws.Range("J80").Value = 1305.8572
$ws.Range("L80").Value = 1305.8572
$ws.Range("N80").Value = -3301.8572
$ws.Range("H83").Value = 1267.625
$ws.Range("J83").Value = 1305.8572
$ws.Range("L83").Value = 6529.286
$ws.Range("N83").Value = -16513.286
$ws.Range("H86").Value = 6283.2104
$ws.Range("I86").Value = 10177.375
$ws.Range("K86").Value = 10177.375
$ws.Range("M86").Value = -9054.375
$ws.Range("H89").Value = 6283.2104
$ws.Range("I89").Value = 10177.375
$ws.Range("K89").Value = 50886.875
$ws.Range("M89").Value = -45270.875
$ws.Range("H107").Value = 8463.294
$ws.Range("I107").Value = 8905.866
$ws.Range("K107").Value = 8905.866
$ws.Range("M107").Value = -6985.866

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 384.22223
$ws.Range("I10").Value = 384.22223
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 384.22223
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -245.22223
$ws.Range("H62").Value = 73504.89
$ws.Range("J62").Value = 123739
$ws.Range("L62").Value = 123739
$ws.Range("N62").Value = -124987
$ws.Range("H65").Value = 73504.89
$ws.Range("J65").Value = 123739
$ws.Range("L65").Value = 618695
$ws.Range("N65").Value = -624935
$ws.Range("H74").Value = 62333
$ws.Range("I74").Value = 40000
$ws.Range("J74").Value = 73499.5
$ws.Range("K74").Value = 40000
$ws.Range("L74").Value = 73499.5
$ws.Range("M74").Value = -39126
$ws.Range("N74").Value = -75247.5
$ws.Range("H77").Value = 62333
$ws.Range("I77").Value = 40000
$ws.Range("J77").Value = 73499.5
$ws.Range("K77").Value = 120000
$ws.Range("L77").Value = 220498.5
$ws.Range("M77").Value = -115632
$ws.Range("N77").Value = -229234.5
$ws.Range("H94").Value = 2227.4285
$ws.Range("I94").Value = 2333.9
$ws.Range("K94").Value = 2333.9
$ws.Range("M94").Value = -1882.9
$ws.Range("H99").Value = 324781
$ws.Range("I99").Value = 735143.1
$ws.Range("J99").Value = 5610.4443
$ws.Range("K99").Value = 735143.1
$ws.Range("L99").Value = 5610.4443
$ws.Range("M99").Value = -733645.1
$ws.Range("N99").Value = -8606.444299999999
$ws.Range("H126").Value = 324781
$ws.Range("I126").Value = 735143.1
$ws.Range("J126").Value = 5610.4443
$ws.Range("K126").Value = 2205429.3
$ws.Range("L126").Value = 16831.3329
$ws.Range("M126").Value = -2202959.3
$ws.Range("N126").Value = -21771.3329

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45507116
$ws.Range("I4").Value = 34871864
$ws.Range("J4").Value = 210353470
$ws.Range("K4").Value = 104615592
$ws.Range("L4").Value = 631060410
$ws.Range("M4").Value = -104615480
$ws.Range("N4").Value = -631060634
$ws.Range("H39").Value = 429.7
$ws.Range("J39").Value = 799.3333
$ws.Range("L39").Value = 2397.9999
$ws.Range("N39").Value = -2985.9999
$ws.Range("H49").Value = 783
$ws.Range("I49").Value = 783
$ws.Range("K49").Value = 2349
$ws.Range("M49").Value = -2193
$ws.Range("H51").Value = 1615.9231
$ws.Range("J51").Value = 3035.75
$ws.Range("L51").Value = 9107.25
$ws.Range("N51").Value = -10027.25
$ws.Range("H55").Value = 8522.76
$ws.Range("J55").Value = 9658.182000000001
$ws.Range("L55").Value = 28974.546
$ws.Range("N55").Value = -29328.546
$ws.Range("H68").Value = 11484.385
$ws.Range("J68").Value = 14160.4
$ws.Range("L68").Value = 42481.2
$ws.Range("N68").Value = -44103.2
$ws.Range("H71").Value = 11484.385
$ws.Range("J71").Value = 14160.4
$ws.Range("L71").Value = 127443.6
$ws.Range("N71").Value = -135555.6
$ws.Range("H86").Value = 1112.4242
$ws.Range("I86").Value = 1279.4
$ws.Range("J86").Value = 1082.6072
$ws.Range("K86").Value = 3838.2
$ws.Range("L86").Value = 3247.8216
$ws.Range("M86").Value = -2652.2
$ws.Range("N86").Value = -5619.821599999999
$ws.Range("H89").Value = 1112.4242
$ws.Range("I89").Value = 1279.4
$ws.Range("J89").Value = 1082.6072
$ws.Range("K89").Value = 11514.6
$ws.Range("L89").Value = 9743.4648
$ws.Range("M89").Value = -5586.6
$ws.Range("N89").Value = -21599.4648
$ws.Range("H113").Value = 3109.5
$ws.Range("J113").Value = 3232.4443
$ws.Range("L113").Value = 9697.332900000001
$ws.Range("N113").Value = -14037.3329

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 24420
$ws.Range("J34").Value = 24420
$ws.Range("L34").Value = 24420
$ws.Range("N34").Value = -24956
$ws.Range("H52").Value = 30333.334
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 30333.334
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").Value = 30333.334
$ws.Range("N52").Value = -30851.334
$ws.Range("H76").Value = 24420
$ws.Range("J76").Value = 24420
$ws.Range("L76").Value = 24420
$ws.Range("N76").Value = -25050
$ws.Range("H79").Value = 24420
$ws.Range("J79").Value = 24420
$ws.Range("L79").Value = 24420
$ws.Range("N79").Value = -26604
$ws.Range("H126").Value = 14281
$ws.Range("I126").Value = 28433
$ws.Range("J126").Value = 9811.947
$ws.Range("K126").Value = 85299
$ws.Range("L126").Value = 29435.841
$ws.Range("M126").Value = -82829
$ws.Range("N126").Value = -34375.841

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3164.1667
$ws.Range("I22").Value = 3348.3635
$ws.Range("J22").Value = 2874.7144
$ws.Range("K22").Value = 3348.3635
$ws.Range("L22").Value = 2874.7144
$ws.Range("M22").Value = -3053.3635
$ws.Range("N22").Value = -3464.7144
$ws.Range("H27").Value = 3164.1667
$ws.Range("I27").Value = 3348.3635
$ws.Range("J27").Value = 2874.7144
$ws.Range("K27").Value = 3348.3635
$ws.Range("L27").Value = 2874.7144
$ws.Range("M27").Value = -3241.3635
$ws.Range("N27").Value = -3088.7144
$ws.Range("H68").Value = 4755.5557
$ws.Range("J68").Value = 6960
$ws.Range("L68").Value = 6960
$ws.Range("N68").Value = -8458
$ws.Range("H71").Value = 4755.5557
$ws.Range("J71").Value = 6960
$ws.Range("L71").Value = 34800
$ws.Range("N71").Value = -42288
$ws.Range("H82").Value = 2559.4736
$ws.Range("I82").Value = 2893.818
$ws.Range("K82").Value = 2893.818
$ws.Range("M82").Value = -2532.818
$ws.Range("H85").Value = 2559.4736
$ws.Range("I85").Value = 2893.818
$ws.Range("K85").Value = 2893.818
$ws.Range("M85").Value = -1645.818
$ws.Range("H132").Value = 599612.1
$ws.Range("I132").Value = 878547.5
$ws.Range("K132").Value = 2635642.5
$ws.Range("M132").Value = -2633112.5
$ws.Range("H136").Value = 8347.608
$ws.Range("I136").Value = 8900.833000000001
$ws.Range("K136").Value = 26702.499
$ws.Range("M136").Value = -24152.499

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 18000
$ws.Range("I18").Value = 18000
$ws.Range("K18").Value = 18000
$ws.Range("M18").Value = -17827
$ws.Range("H93").Value = 58694.5
$ws.Range("J93").Value = 58694.5
$ws.Range("L93").Value = 58694.5
$ws.Range("N93").Value = -63686.5
$ws.Range("H132").Value = 24282.36
$ws.Range("I132").Value = 38120.645
$ws.Range("K132").Value = 114361.935
$ws.Range("M132").Value = -111831.935
